$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9 (shifts existing row 9 and below down by 1)
$ws.Rows.Item(9).Insert()

# Update B8 value
$ws.Cells.Item(8, 2).Value = 21527.67

# Populate the new row 9 with data
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 17965.81
$ws.Cells.Item(9, 3).Value = 7
$ws.Cells.Item(9, 4).Value = 2025
$ws.Cells.Item(9, 5).Value = "07/2025"
